$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: 'Datos actualizados a 15 de Junio de 2020 a las 00:05' -> 'Datos actualizados a 15 de Junio de 2020 a las 01:22'
$ws.Range("A1").Value = "Datos actualizados a 15 de Junio de 2020 a las 01:22"

# Row 4: 'Estados Unidos' -> 'Estados Unidos'
$ws.Range("B4").Value = 2161366
$ws.Range("C4").Value = 19142
$ws.Range("D4").Value = 867074
$ws.Range("E4").Value = 1176443
$ws.Range("G4").Value = 322
$ws.Range("H4").Value = 117849

# Row 5: 'Brasil' -> 'Brasil'
$ws.Range("B5").Value = 867882
$ws.Range("C5").Value = 17086
$ws.Range("E5").Value = 386981
$ws.Range("G5").Value = 598
$ws.Range("H5").Value = 43389

# Row 20: 'Canada' -> 'Canada'
$ws.Range("B20").Value = 98787
$ws.Range("C20").Value = 377
$ws.Range("D20").Value = 60272
$ws.Range("E20").Value = 30369

# Row 28: 'Paises Bajos' -> 'Colombia'
$ws.Range("A28").Value = "Colombia"
$ws.Range("B28").Value = 50939
$ws.Range("C28").Value = 2193
$ws.Range("D28").Value = 19822
$ws.Range("E28").Value = 29450
$ws.Range("G28").Value = 75
$ws.Range("H28").Value = 1667

# Row 29: 'Colombia' -> 'Paises Bajos'
$ws.Range("A29").Value = "Paises Bajos"
$ws.Range("B29").Value = 48783
$ws.Range("C29").Value = 143
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 6059

# Row 37: 'Ucrania' -> 'Argentina'
$ws.Range("A37").Value = "Argentina"
$ws.Range("B37").Value = 31577
$ws.Range("C37").Value = 1282
$ws.Range("D37").Value = 9564
$ws.Range("E37").Value = 21180
$ws.Range("G37").Value = 18
$ws.Range("H37").Value = 833

# Row 38: 'Suiza' -> 'Ucrania'
$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 31154
$ws.Range("C38").Value = 648
$ws.Range("D38").Value = 14082
$ws.Range("E38").Value = 16183
$ws.Range("G38").Value = 9
$ws.Range("H38").Value = 889

# Row 39: 'Argentina' -> 'Suiza'
$ws.Range("A39").Value = "Suiza"
$ws.Range("B39").Value = 31117
$ws.Range("C39").Value = 23
$ws.Range("D39").Value = 28800
$ws.Range("E39").Value = 379
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 1938

# Row 52: 'Japon' -> 'Japon'
$ws.Range("B52").Value = 17429
$ws.Range("C52").Value = 47
$ws.Range("D52").Value = 15643
$ws.Range("E52").Value = 861
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 925

# Row 55: 'Nigeria' -> 'Nigeria'
$ws.Range("B55").Value = 16085
$ws.Range("C55").Value = 403
$ws.Range("D55").Value = 5220
$ws.Range("E55").Value = 10445
$ws.Range("G55").Value = 13
$ws.Range("H55").Value = 420

# Row 63: 'Chequia' -> 'Chequia'
$ws.Range("B63").Value = 10024
$ws.Range("C63").Value = 33
$ws.Range("D63").Value = 7226
$ws.Range("E63").Value = 2469

# Row 72: 'Finlandia' -> 'Sudan'
$ws.Range("A72").Value = "Sudan"
$ws.Range("B72").Value = 7220
$ws.Range("C72").Value = 213
$ws.Range("D72").Value = 2610
$ws.Range("E72").Value = 4151
$ws.Range("G72").Value = 12
$ws.Range("H72").Value = 459

# Row 73: 'Sudan' -> 'Finlandia'
$ws.Range("A73").Value = "Finlandia"
$ws.Range("B73").Value = 7104
$ws.Range("C73").Value = 17
$ws.Range("D73").Value = 6200
$ws.Range("E73").Value = 578
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 326

# Row 80: 'Guinea' -> 'Guinea'
$ws.Range("B80").Value = 4532
$ws.Range("C80").Value = 48
$ws.Range("D80").Value = 3234
$ws.Range("E80").Value = 1273

# Row 96: 'Kirguistan' -> 'Mayotte'
$ws.Range("A96").Value = "Mayotte"
$ws.Range("B96").Value = 2298
$ws.Range("C96").Value = 16
$ws.Range("D96").Value = 1790
$ws.Range("E96").Value = 479
$ws.Range("G96").Value = 1
$ws.Range("H96").Value = 29

# Row 97: 'Mayotte' -> 'Kirguistan'
$ws.Range("A97").Value = "Kirguistan"
$ws.Range("B97").Value = 2285
$ws.Range("C97").Value = 78
$ws.Range("D97").Value = 1791
$ws.Range("E97").Value = 467
$ws.Range("H97").Value = 27

# Row 140: 'Santo Tome y Principe' -> 'Santo Tome y Principe'
$ws.Range("B140").Value = 661
$ws.Range("C140").Value = 2
$ws.Range("D140").Value = 177
$ws.Range("E140").Value = 472

# Row 170: 'Guyana' -> 'Guyana'
$ws.Range("D170").Value = 99
$ws.Range("E170").Value = 48

# Row 175: 'Trinidad yTobago' -> 'Trinidad yTobago'
$ws.Range("B175").Value = 123
$ws.Range("C175").Value = 6
$ws.Range("E175").Value = 6

# Row 206: 'Groenlandia' -> 'Islas Malvinas'
$ws.Range("A206").Value = "Islas Malvinas"

# Row 207: 'Islas Malvinas' -> 'Groenlandia'
$ws.Range("A207").Value = "Groenlandia"

# Row 208: 'Islas Turcas y Caicos' -> 'Santa Sede'
$ws.Range("A208").Value = "Santa Sede"
$ws.Range("D208").Value = 12
$ws.Range("H208").Value = 0

# Row 209: 'Santa Sede' -> 'Islas Turcas y Caicos'
$ws.Range("A209").Value = "Islas Turcas y Caicos"
$ws.Range("D209").Value = 11
$ws.Range("H209").Value = 1
